$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 227, pushing the existing rows 227-249 down to 229-251.
$ws.Rows("227:228").Insert()

# Populate new row 227
$ws.Range("A227").Value = 1
$ws.Range("B227").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C227").Value = 'Arica y Parinacota'
$ws.Range("D227").Value = 45132
$ws.Range("E227").Value = 15
$ws.Range("F227").Value = 100114001
$ws.Range("G227").Value = 'Papa'
$ws.Range("H227").Value = 'Cardinal'
$ws.Range("I227").Value = '1a (cosecha)'
$ws.Range("J227").Value = 1000
$ws.Range("K227").Value = 21000
$ws.Range("L227").Value = 22000
$ws.Range("M227").Value = 21500
$ws.Range("N227").Value = '$/saco 25 kilos'
$ws.Range("O227").Value = 'Región Metropolitana'
$ws.Range("P227").Value = 860
$ws.Range("Q227").Value = 25
$ws.Range("R227").Value = 'Hortaliza'

# Populate new row 228
$ws.Range("A228").Value = 1
$ws.Range("B228").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C228").Value = 'Arica y Parinacota'
$ws.Range("D228").Value = 45132
$ws.Range("E228").Value = 15
$ws.Range("F228").Value = 100114001
$ws.Range("G228").Value = 'Papa'
$ws.Range("H228").Value = 'Patagonia'
$ws.Range("I228").Value = '1a (cosecha)'
$ws.Range("J228").Value = 1000
$ws.Range("K228").Value = 20000
$ws.Range("L228").Value = 21000
$ws.Range("M228").Value = 20500
$ws.Range("N228").Value = '$/saco 25 kilos'
$ws.Range("O228").Value = 'Región de Los Lagos'
$ws.Range("P228").Value = 820
$ws.Range("Q228").Value = 25
$ws.Range("R228").Value = 'Hortaliza'
